# Update JS-SPA Self-Evaluation Protocol worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JS-SPA-Self-Evaluation-Protocol")

# "Days Commit in GitHub" row (row 8) comment stays the same text, nothing to
# change there - shared-string index shift only happens on save, not a value
# edit we need to perform explicitly.

# Admin Home Screen row: change "Yes Half" to "Yes" and add comment about paging
$ws.Range("C34").Value = "Yes"
$ws.Range("E34").Value = "paging-ът малко се счупи"

# Admin Approve/Reject/Edit/Delete Ad rows: set answer to "Yes"
$ws.Range("C35").Value = "Yes"
$ws.Range("C36").Value = "Yes"
$ws.Range("C37").Value = "Yes"
$ws.Range("C38").Value = "Yes"

# Row 11 "Web Design" comment -> add new comment
$ws.Range("E11").Value = "дизайнът е responsive"

# Update the active selection to F11
$ws.Range("F11").Select()
